$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New answer text for the "Meat" question group (rows 2-10, column B) -
# the plain product names are turned into hyperlinks to the Sysco product
# categories page.
$meatText = "We have multiple kinds of Meat. Please select your choice. 1.<a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Beef</a> 2. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Pork</a> 3. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Chicken</a>"

# New answer text for the "Dairy" question group (rows 11-19, column B) -
# same treatment for the dairy product names.
$dairyText = "Among Dairy products we have 1. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Milk</a> 2. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Yogurt</a> 3. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Greek Yogurt</a> 4. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Cheese</a> 5. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Milk Powder</a> 6. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Butter</a> 7. <a href = 'https://www.sysco.com/Products/Products/Product-Categories.html',target='_blank'>Ice Cream</a>"

$ws.Range("B2:B10").Value = $meatText
$ws.Range("B11:B19").Value = $dairyText

# Leave the selection where the last edit happened, matching the saved
# workbook's view state.
$ws.Range("B19").Select() | Out-Null
